{"js": "// Update the worksheet date and each two-digit-by-two-digit multiplication\n// answer cell to the new values from the latest generated output.\nconst replacements = [\n  [\"2025-10-16 Thursday\", \"2025-10-17 Friday\"],\n  [\"62\u00d772=4464\", \"26\u00d734=884\"],\n  [\"24\u00d730=720\", \"76\u00d757=4332\"],\n  [\"12\u00d715=180\", \"54\u00d777=4158\"],\n  [\"21\u00d756=1176\", \"85\u00d790=7650\"],\n  [\"48\u00d743=2064\", \"29\u00d760=1740\"],\n  [\"30\u00d751=1530\", \"73\u00d749=3577\"],\n  [\"33\u00d742=1386\", \"23\u00d738=874\"],\n  [\"65\u00d789=5785\", \"67\u00d789=5963\"],\n  [\"15\u00d765=975\", \"18\u00d797=1746\"],\n  [\"80\u00d754=4320\", \"17\u00d768=1156\"],\n  [\"46\u00d724=1104\", \"91\u00d772=6552\"],\n  [\"30\u00d741=1230\", \"50\u00d748=2400\"],\n  [\"63\u00d768=4284\", \"54\u00d725=1350\"],\n  [\"81\u00d725=2025\", \"37\u00d735=1295\"],\n  [\"27\u00d714=378\", \"50\u00d746=2300\"],\n  [\"52\u00d724=1248\", \"28\u00d789=2492\"],\n  [\"63\u00d738=2394\", \"41\u00d756=2296\"],\n  [\"82\u00d789=7298\", \"78\u00d779=6162\"],\n  [\"63\u00d739=2457\", \"26\u00d718=468\"],\n  [\"42\u00d746=1932\", \"53\u00d791=4823\"],\n  [\"68\u00d767=4556\", \"19\u00d731=589\"],\n  [\"77\u00d799=7623\", \"27\u00d728=756\"],\n  [\"96\u00d793=8928\", \"97\u00d782=7954\"],\n  [\"52\u00d789=4628\", \"11\u00d780=880\"],\n  [\"80\u00d712=960\", \"83\u00d745=3735\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and each two-digit-by-two-digit multiplication\n# answer cell to the new values from the latest generated output.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n    @(\"2025-10-16 Thursday\", \"2025-10-17 Friday\"),\n    @(\"62\u00d772=4464\", \"26\u00d734=884\"),\n    @(\"24\u00d730=720\", \"76\u00d757=4332\"),\n    @(\"12\u00d715=180\", \"54\u00d777=4158\"),\n    @(\"21\u00d756=1176\", \"85\u00d790=7650\"),\n    @(\"48\u00d743=2064\", \"29\u00d760=1740\"),\n    @(\"30\u00d751=1530\", \"73\u00d749=3577\"),\n    @(\"33\u00d742=1386\", \"23\u00d738=874\"),\n    @(\"65\u00d789=5785\", \"67\u00d789=5963\"),\n    @(\"15\u00d765=975\", \"18\u00d797=1746\"),\n    @(\"80\u00d754=4320\", \"17\u00d768=1156\"),\n    @(\"46\u00d724=1104\", \"91\u00d772=6552\"),\n    @(\"30\u00d741=1230\", \"50\u00d748=2400\"),\n    @(\"63\u00d768=4284\", \"54\u00d725=1350\"),\n    @(\"81\u00d725=2025\", \"37\u00d735=1295\"),\n    @(\"27\u00d714=378\", \"50\u00d746=2300\"),\n    @(\"52\u00d724=1248\", \"28\u00d789=2492\"),\n    @(\"63\u00d738=2394\", \"41\u00d756=2296\"),\n    @(\"82\u00d789=7298\", \"78\u00d779=6162\"),\n    @(\"63\u00d739=2457\", \"26\u00d718=468\"),\n    @(\"42\u00d746=1932\", \"53\u00d791=4823\"),\n    @(\"68\u00d767=4556\", \"19\u00d731=589\"),\n    @(\"77\u00d799=7623\", \"27\u00d728=756\"),\n    @(\"96\u00d793=8928\", \"97\u00d782=7954\"),\n    @(\"52\u00d789=4628\", \"11\u00d780=880\"),\n    @(\"80\u00d712=960\", \"83\u00d745=3735\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n}\n"}
